$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.251.35'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '2.996.29'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '501.90'
$ws.Range('E5').Value = '  -4.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.95'
$ws.Range('E6').Value = '  -3.99%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -3.80%  '
$ws.Range('E9').Value = '  -5.18%  '
$ws.Range('E10').Value = '  -4.56%  '
$ws.Range('E11').Value = '  -3.78%  '
$ws.Range('D12').Value = '3.506.80'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.01'
$ws.Range('E14').Value = '  -4.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000160'
$ws.Range('E15').Value = '  -6.34%  '
$ws.Range('D16').Value = '57.258.18'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.08'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('D18').Value = '2.992.05'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E19').Value = '  -4.18%  '
$ws.Range('E20').Value = '  -3.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '319.79'
$ws.Range('E21').Value = '  -5.89%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.75'
$ws.Range('E23').Value = '  +1.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.491'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.20'
$ws.Range('E25').Value = '  -2.76%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E28').Value = '  -8.88%  '
$ws.Range('E29').Value = '  -5.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.08'
$ws.Range('E30').Value = '  -3.98%  '
$ws.Range('E31').Value = '  -4.25%  '
$ws.Range('E32').Value = '  -7.01%  '
$ws.Range('E33').Value = '  -4.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '154.99'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.54'
$ws.Range('E35').Value = '  -4.80%  '
$ws.Range('E36').Value = '  -3.95%  '
$ws.Range('E37').Value = '  -7.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.40'
$ws.Range('E38').Value = '  -7.49%  '
$ws.Range('E39').Value = '  -5.83%  '
$ws.Range('B40').Value = 'RenzoRestakedETH'
$ws.Range('C40').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D40').Value = '3.026.29'
$ws.Range('E40').Value = '  -2.26%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.87'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').Value = '  -4.68%  '
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('D45').Value = '2.182.35'
$ws.Range('E45').Value = '  -6.51%  '
$ws.Range('E46').Value = '  -6.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.94'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.936'
$ws.Range('E48').Value = '  -9.08%  '
$ws.Range('E49').Value = '  -4.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.16'
$ws.Range('E50').Value = '  -5.37%  '
$ws.Range('E51').Value = '  -12.46%  '
